# Burndown Chart Caliber Mobile - sprint data update.
# "Planned" slips for one more day (13 -> 5 on 8/9) and the "Actual" burndown
# series is extended through 8/12 with the newly-reported completions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Planned (column B) - day 8/9 (row 15) actually finished at 5, not 13.
$ws.Range("B15").Value = 5

# Actual (column C) - fill in the days that previously had no reported value.
$ws.Range("C12").Value = 11
$ws.Range("C13").Value = 11
$ws.Range("C14").Value = 11
$ws.Range("C15").Value = 11
$ws.Range("C16").Value = 3
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 3

# Refresh both burndown charts (Sheet1 and the Sheet2 skeleton/presentation
# copy) so they pick up the new Planned/Actual values from the worksheet.
foreach ($sheet in $wb.Worksheets) {
    $chartObjects = $sheet.ChartObjects()
    for ($i = 1; $i -le $chartObjects.Count; $i++) {
        $chartObjects.Item($i).Chart.Refresh()
    }
}

# Leave the cursor where the author last left it when saving.
$ws.Range("O17").Select()
